$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row tweak
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "MasterSheet RowNo."

# ---------------------------------------------------------------------------
# 2. Fill in the previously-missing TotalConfirmedNewCases (G) / TotalNewDeaths
#    (I) columns for the existing rows 2-9
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = 1
$ws.Range("I2").Value = 0

$ws.Range("G3").Value = 0
$ws.Range("I3").Value = 0

$ws.Range("G4").Value = 0
$ws.Range("I4").Value = 0

$ws.Range("G5").Value = 0
$ws.Range("I5").Value = 0

$ws.Range("G6").Value = 6
$ws.Range("I6").Value = 0

$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 0

$ws.Range("G8").Value = 0
$ws.Range("I8").Value = 0

$ws.Range("G9").Value = 0
$ws.Range("I9").Value = 0

# ---------------------------------------------------------------------------
# 3. Brand new row of data (row 10)
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 71
$ws.Range("B10").Value = 82
$ws.Range("C10").Value = "LATIN AMER. & CARIB    "
$ws.Range("D10").Value = 43921
$ws.Range("E10").Value = "Grenada"
$ws.Range("F10").Value = 9
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "Local transmission"
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 5361

# ---------------------------------------------------------------------------
# 4. Column widths - A:O all 27 "characters" wide
# ---------------------------------------------------------------------------
$ws.Range("A1:O1").EntireColumn.ColumnWidth = 26.16666666666667

# ---------------------------------------------------------------------------
# 5. Formatting - center every cell (horizontally + vertically), and give the
#    Date column its own custom number format. Using a template cell + Copy /
#    PasteSpecial(Formats) applies a single, fully-formed style per group in
#    one shot instead of incrementally building up alignment attribute by
#    attribute (which would otherwise leave behind unused intermediate
#    styles).
# ---------------------------------------------------------------------------
$tmplNormal = $ws.Cells.Item(20, 1)
$tmplNormal.HorizontalAlignment = -4108
$tmplNormal.VerticalAlignment = -4108

$tmplDate = $ws.Cells.Item(20, 2)
$tmplDate.NumberFormat = "yyyy-mm-dd;"
$tmplDate.HorizontalAlignment = -4108
$tmplDate.VerticalAlignment = -4108

$tmplNormal.Copy()
$ws.Range("A1:O10").PasteSpecial(-4122)

$tmplDate.Copy()
$ws.Range("D1:D10").PasteSpecial(-4122)

# Clean up the scratch template cells used as the copy source above.
$ws.Rows.Item(20).Delete()
